# Insert a new data row before existing row 112 (shifting rows 112..208 down to 113..209)
# and populate the new row 112 with a new record for "Pepino ensalada".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 112, pushing everything below it down by one row.
$ws.Rows.Item(112).Insert()

# Fill in the values for the newly inserted row 112.
$ws.Cells.Item(112, 1).Value = 11
$ws.Cells.Item(112, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(112, 3).Value = "Bíobío"
$ws.Cells.Item(112, 4).Value = 45062
$ws.Cells.Item(112, 5).Value = 8
$ws.Cells.Item(112, 6).Value = 100112043
$ws.Cells.Item(112, 7).Value = "Pepino ensalada"
$ws.Cells.Item(112, 8).Value = "Sin especificar"
$ws.Cells.Item(112, 9).Value = "Primera"
$ws.Cells.Item(112, 10).Value = 100
$ws.Cells.Item(112, 11).Value = 12000
$ws.Cells.Item(112, 12).Value = 13000
$ws.Cells.Item(112, 13).Value = 12500
$ws.Cells.Item(112, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(112, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(112, 16).Value = 208
$ws.Cells.Item(112, 17).Value = 60
$ws.Cells.Item(112, 18).Value = "Hortaliza"
